$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three rows that were removed from the data table. They must be
# removed from the bottom-most row upward so the remaining row numbers do not
# shift under us while deleting.
#   Row 24: PROYECTO DE SISTEMAS ROBUSTOS, PARALELOS Y DISTRIBUIDOS
#   Row 23: PROYECTO DE GESTION DE LA TECNOLOGIA DE INFORMACION
#   Row 9 : COMPUTO FLEXIBLE (SOFTCOMPUTING)
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(9).Delete()
